$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Row 61: 2022-05-26
$ws.Range("A61").Value = 44707
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 327417
$ws.Range("D61").Formula = "=D60+F60"
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 0
$ws.Range("A61").NumberFormat = "yyyy\-mm\-dd;@"

# Row 62: 2022-05-27
$ws.Range("A62").Value = 44708
$ws.Range("B62").Value = 0
$ws.Range("C62").Value = 327429
$ws.Range("D62").Value = 6347
$ws.Range("E62").Value = 12
$ws.Range("F62").Value = 0
$ws.Range("A62").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("F62").Select()
